$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 2017.2632
$ws.Cells.Item(28, 9).Value = 781.25
$ws.Cells.Item(28, 10).Value = 4136.143
$ws.Cells.Item(28, 11).Value = 781.25
$ws.Cells.Item(28, 12).Value = 4136.143
$ws.Cells.Item(28, 13).Value = -296.25
$ws.Cells.Item(28, 14).Value = -5106.143

$ws.Cells.Item(43, 8).Value = 1998
$ws.Cells.Item(43, 9).Value = 997
$ws.Cells.Item(43, 10).Value = 4000
$ws.Cells.Item(43, 11).Value = 997
$ws.Cells.Item(43, 12).Value = 4000
$ws.Cells.Item(43, 13).Value = -928
$ws.Cells.Item(43, 14).Value = -4138

$ws.Cells.Item(80, 8).Value = 580
$ws.Cells.Item(80, 9).Value = 575
$ws.Cells.Item(80, 10).Value = 600
$ws.Cells.Item(80, 11).Value = 1725
$ws.Cells.Item(80, 12).Value = 1800
$ws.Cells.Item(80, 13).Value = -727
$ws.Cells.Item(80, 14).Value = -3796

$ws.Cells.Item(83, 8).Value = 580
$ws.Cells.Item(83, 9).Value = 575
$ws.Cells.Item(83, 10).Value = 600
$ws.Cells.Item(83, 11).Value = 5175
$ws.Cells.Item(83, 12).Value = 5400
$ws.Cells.Item(83, 13).Value = -183
$ws.Cells.Item(83, 14).Value = -15384

$ws.Cells.Item(88, 8).Value = 1750
$ws.Cells.Item(88, 9).Value = 1500
$ws.Cells.Item(88, 10).Value = 2000
$ws.Cells.Item(88, 11).Value = 1500
$ws.Cells.Item(88, 12).Value = 2000
$ws.Cells.Item(88, 13).Value = -1094
$ws.Cells.Item(88, 14).Value = -2812

$ws.Cells.Item(91, 8).Value = 1750
$ws.Cells.Item(91, 9).Value = 1500
$ws.Cells.Item(91, 10).Value = 2000
$ws.Cells.Item(91, 11).Value = 1500
$ws.Cells.Item(91, 12).Value = 2000
$ws.Cells.Item(91, 13).Value = -96
$ws.Cells.Item(91, 14).Value = -4808

$ws.Cells.Item(92, 8).Value = 154
$ws.Cells.Item(92, 9).Value = 154
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 154
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).Value = 1094

$ws.Cells.Item(100, 8).Value = 3210
$ws.Cells.Item(100, 9).Value = 3210
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 11).Value = 3210
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 13).Value = -2669
$ws.Cells.Item(100, 14).Value = ""

$ws.Cells.Item(116, 8).Value = 2757
$ws.Cells.Item(116, 9).Value = 2667.25
$ws.Cells.Item(116, 10).Value = 2996.3333
$ws.Cells.Item(116, 11).Value = 2667.25
$ws.Cells.Item(116, 12).Value = 2996.3333
$ws.Cells.Item(116, 13).Value = 774.75
$ws.Cells.Item(116, 14).Value = -9880.3333

$ws.Cells.Item(137, 8).Value = 1672.25
$ws.Cells.Item(137, 9).Value = 1196.8572
$ws.Cells.Item(137, 10).Value = 5000
$ws.Cells.Item(137, 11).Value = 3590.5716
$ws.Cells.Item(137, 12).Value = 15000
$ws.Cells.Item(137, 13).Value = -1040.5716
$ws.Cells.Item(137, 14).Value = -20100

$ws.Cells.Item(138, 8).Value = 3013.125
$ws.Cells.Item(138, 9).Value = 1213.75
$ws.Cells.Item(138, 10).Value = 4812.5
$ws.Cells.Item(138, 11).Value = 3641.25
$ws.Cells.Item(138, 12).Value = 14437.5
$ws.Cells.Item(138, 13).Value = 1498.75
$ws.Cells.Item(138, 14).Value = -24717.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 141.44444
$ws.Cells.Item(5, 9).Value = 112.6
$ws.Cells.Item(5, 10).Value = 177.5
$ws.Cells.Item(5, 11).Value = 112.6
$ws.Cells.Item(5, 12).Value = 177.5
$ws.Cells.Item(5, 13).Value = -0.5999999999999943
$ws.Cells.Item(5, 14).Value = -401.5

$ws.Cells.Item(18, 8).Value = 8888
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 8888
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 8888
$ws.Cells.Item(18, 14).Value = -9532

$ws.Cells.Item(32, 8).Value = 3056.3794
$ws.Cells.Item(32, 9).Value = 3056.3794
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 3056.3794
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = -2769.3794

$ws.Cells.Item(45, 8).Value = 2197.7273
$ws.Cells.Item(45, 9).Value = 1393.3334
$ws.Cells.Item(45, 10).Value = 3163
$ws.Cells.Item(45, 11).Value = 1393.3334
$ws.Cells.Item(45, 12).Value = 3163
$ws.Cells.Item(45, 13).Value = -1016.3334
$ws.Cells.Item(45, 14).Value = -3917

$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).Value = ""

$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).Value = ""

$ws.Cells.Item(97, 8).Value = 451.85715
$ws.Cells.Item(97, 9).Value = 451.85715
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 451.85715
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = 44.14285000000001

$ws.Cells.Item(110, 8).Value = 1499.3334
$ws.Cells.Item(110, 9).Value = 1283
$ws.Cells.Item(110, 10).Value = 2256.5
$ws.Cells.Item(110, 11).Value = 1283
$ws.Cells.Item(110, 12).Value = 2256.5
$ws.Cells.Item(110, 13).Value = 762
$ws.Cells.Item(110, 14).Value = -6346.5

$ws.Cells.Item(122, 8).Value = 555.1111
$ws.Cells.Item(122, 9).Value = 555.1111
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 1665.3333
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = 784.6667000000002

$ws.Cells.Item(132, 8).Value = 5726.8184
$ws.Cells.Item(132, 9).Value = 5726.8184
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 17180.4552
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -14650.4552

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 141.44444
$ws.Cells.Item(4, 9).Value = 112.6
$ws.Cells.Item(4, 10).Value = 177.5
$ws.Cells.Item(4, 11).Value = 112.6
$ws.Cells.Item(4, 12).Value = 177.5
$ws.Cells.Item(4, 13).Value = 2.400000000000006
$ws.Cells.Item(4, 14).Value = -407.5

$ws.Cells.Item(20, 8).Value = 2619.8572
$ws.Cells.Item(20, 9).Value = 2031
$ws.Cells.Item(20, 10).Value = 3061.5
$ws.Cells.Item(20, 11).Value = 2031
$ws.Cells.Item(20, 12).Value = 3061.5
$ws.Cells.Item(20, 13).Value = -1784
$ws.Cells.Item(20, 14).Value = -3555.5

$ws.Cells.Item(22, 8).Value = 1789.6666
$ws.Cells.Item(22, 9).Value = 1789.6666
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 1789.6666
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -1616.6666

$ws.Cells.Item(94, 8).Value = 256.42856
$ws.Cells.Item(94, 9).Value = 256.42856
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 256.42856
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).Value = 194.57144
$ws.Cells.Item(94, 14).Value = ""

$ws.Cells.Item(107, 8).Value = 3666.7
$ws.Cells.Item(107, 9).Value = 2549.0588
$ws.Cells.Item(107, 10).Value = 10000
$ws.Cells.Item(107, 11).Value = 2549.0588
$ws.Cells.Item(107, 12).Value = 10000
$ws.Cells.Item(107, 13).Value = -629.0587999999998
$ws.Cells.Item(107, 14).Value = -13840

$ws.Cells.Item(134, 8).Value = 3974
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 3974
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 11922
$ws.Cells.Item(134, 14).Value = -16992

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 5999.75
$ws.Cells.Item(86, 9).Value = 4666.3335
$ws.Cells.Item(86, 10).Value = 10000
$ws.Cells.Item(86, 11).Value = 4666.3335
$ws.Cells.Item(86, 12).Value = 10000
$ws.Cells.Item(86, 13).Value = -3543.3335
$ws.Cells.Item(86, 14).Value = -12246

$ws.Cells.Item(89, 8).Value = 5999.75
$ws.Cells.Item(89, 9).Value = 4666.3335
$ws.Cells.Item(89, 10).Value = 10000
$ws.Cells.Item(89, 11).Value = 23331.6675
$ws.Cells.Item(89, 12).Value = 50000
$ws.Cells.Item(89, 13).Value = -17715.6675
$ws.Cells.Item(89, 14).Value = -61232

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(22, 8).Value = 294.66666
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 294.66666
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 883.9999799999999
$ws.Cells.Item(22, 14).Value = -1221.99998

$ws.Cells.Item(27, 8).Value = 294.66666
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 294.66666
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 883.9999799999999
$ws.Cells.Item(27, 14).Value = -1087.99998

$ws.Cells.Item(33, 8).Value = 292.4
$ws.Cells.Item(33, 9).Value = 31.5
$ws.Cells.Item(33, 10).Value = 466.33334
$ws.Cells.Item(33, 11).Value = 189
$ws.Cells.Item(33, 12).Value = 2798.00004
$ws.Cells.Item(33, 13).Value = 94
$ws.Cells.Item(33, 14).Value = -3364.00004

$ws.Cells.Item(132, 8).Value = 4768.3335
$ws.Cells.Item(132, 9).Value = 4498.75
$ws.Cells.Item(132, 10).Value = 4984
$ws.Cells.Item(132, 11).Value = 40488.75
$ws.Cells.Item(132, 12).Value = 44856
$ws.Cells.Item(132, 13).Value = -37958.75
$ws.Cells.Item(132, 14).Value = -49916

$ws.Cells.Item(140, 8).Value = 1614.8
$ws.Cells.Item(140, 9).Value = 768.5
$ws.Cells.Item(140, 10).Value = 5000
$ws.Cells.Item(140, 11).Value = 2305.5
$ws.Cells.Item(140, 12).Value = 15000
$ws.Cells.Item(140, 13).Value = 2874.5
$ws.Cells.Item(140, 14).Value = -25360

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 14).Value = ""

$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 14).Value = ""

$ws.Cells.Item(99, 8).Value = 7985
$ws.Cells.Item(99, 9).Value = 7985
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 7985
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = -5739

$ws.Cells.Item(102, 8).Value = 3266.077
$ws.Cells.Item(102, 9).Value = 3121.5833
$ws.Cells.Item(102, 10).Value = 5000
$ws.Cells.Item(102, 11).Value = 3121.5833
$ws.Cells.Item(102, 12).Value = 5000
$ws.Cells.Item(102, 13).Value = -1499.5833
$ws.Cells.Item(102, 14).Value = -8244

$ws.Cells.Item(113, 8).Value = 6138.1
$ws.Cells.Item(113, 9).Value = 3563.6667
$ws.Cells.Item(113, 10).Value = 9999.75
$ws.Cells.Item(113, 11).Value = 3563.6667
$ws.Cells.Item(113, 12).Value = 9999.75
$ws.Cells.Item(113, 13).Value = -1393.6667
$ws.Cells.Item(113, 14).Value = -14339.75

$ws.Cells.Item(122, 8).Value = 4131
$ws.Cells.Item(122, 9).Value = 4131
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 12393
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -9943

$ws.Cells.Item(132, 8).Value = 2784.4285
$ws.Cells.Item(132, 9).Value = 2598.2
$ws.Cells.Item(132, 10).Value = 3250
$ws.Cells.Item(132, 11).Value = 7794.599999999999
$ws.Cells.Item(132, 12).Value = 9750
$ws.Cells.Item(132, 13).Value = -5264.599999999999
$ws.Cells.Item(132, 14).Value = -14810

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 8747.5
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 8747.5
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 8747.5
$ws.Cells.Item(46, 13).Value = ""
$ws.Cells.Item(46, 14).Value = -9123.5

$ws.Cells.Item(68, 8).Value = 3417.75
$ws.Cells.Item(68, 9).Value = 3417.75
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 3417.75
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = -2668.75
$ws.Cells.Item(68, 14).Value = ""

$ws.Cells.Item(71, 8).Value = 3417.75
$ws.Cells.Item(71, 9).Value = 3417.75
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 17088.75
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = -13344.75
$ws.Cells.Item(71, 14).Value = ""

$ws.Cells.Item(93, 8).Value = 1283.3334
$ws.Cells.Item(93, 9).Value = 1283.3334
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 1283.3334
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = -35.33339999999998

$ws.Cells.Item(132, 8).Value = 3912.7144
$ws.Cells.Item(132, 9).Value = 3221
$ws.Cells.Item(132, 10).Value = 4835
$ws.Cells.Item(132, 11).Value = 9663
$ws.Cells.Item(132, 12).Value = 14505
$ws.Cells.Item(132, 13).Value = -7133
$ws.Cells.Item(132, 14).Value = -19565

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = ""
$ws.Cells.Item(32, 14).Value = ""

$ws.Cells.Item(107, 8).Value = 1766
$ws.Cells.Item(107, 9).Value = 2125.9092
$ws.Cells.Item(107, 10).Value = 446.33334
$ws.Cells.Item(107, 11).Value = 6377.7276
$ws.Cells.Item(107, 12).Value = 1339.00002
$ws.Cells.Item(107, 13).Value = -4457.7276
$ws.Cells.Item(107, 14).Value = -5179.000019999999

$ws.Cells.Item(132, 8).Value = 2302.375
$ws.Cells.Item(132, 9).Value = 2070
$ws.Cells.Item(132, 10).Value = 2999.5
$ws.Cells.Item(132, 11).Value = 6210
$ws.Cells.Item(132, 12).Value = 8998.5
$ws.Cells.Item(132, 13).Value = -3680
$ws.Cells.Item(132, 14).Value = -14058.5

$ws.Cells.Item(136, 8).Value = 3631.2083
$ws.Cells.Item(136, 9).Value = 2264.5386
$ws.Cells.Item(136, 10).Value = 5246.364
$ws.Cells.Item(136, 11).Value = 6793.6158
$ws.Cells.Item(136, 12).Value = 15739.092
$ws.Cells.Item(136, 13).Value = -4243.6158
$ws.Cells.Item(136, 14).Value = -20839.092
